$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.002.18'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '2.013.22'
$ws.Range("E3").Value = '  -1.82%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.49'
$ws.Range("E5").Value = '  -1.70%  '

$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.75'
$ws.Range("E8").Value = '  -3.87%  '

$ws.Range("E9").Value = '  -1.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0789'
$ws.Range("E10").Value = '  +2.69%  '

$ws.Range("E11").Value = '  -3.17%  '

$ws.Range("D12").Value = '2.310.11'
$ws.Range("E12").Value = '  -1.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.25'
$ws.Range("E13").Value = '  -3.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.31'
$ws.Range("E14").Value = '  -1.08%  '

$ws.Range("E15").Value = '  -1.95%  '

$ws.Range("D17").Value = '2.011.19'
$ws.Range("E17").Value = '  -1.69%  '

$ws.Range("D18").Value = '36.923.33'
$ws.Range("E18").Value = '  -1.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.09'
$ws.Range("E19").Value = '  +1.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.77'
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("E22").Value = '  -1.16%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +1.80%  '

$ws.Range("E25").Value = '  -5.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.03'
$ws.Range("E26").Value = '  -2.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.19'
$ws.Range("E27").Value = '  -3.31%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.126'
$ws.Range("E28").Value = '  -3.41%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.36'
$ws.Range("E29").Value = '  +1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.69'
$ws.Range("E30").Value = '  -2.26%  '

$ws.Range("E31").Value = '  -3.45%  '

$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0615'
$ws.Range("E33").Value = '  -1.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.43'
$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("E35").Value = '  -5.68%  '

$ws.Range("E36").Value = '  +1.97%  '

$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  -4.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("E39").Value = '  +2.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0217'
$ws.Range("E40").Value = '  -3.76%  '

$ws.Range("D41").Value = '1.475.19'
$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '95.13'
$ws.Range("E42").Value = '  -3.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.51'
$ws.Range("E43").Value = '  -0.29%  '

$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("E45").Value = '  -5.16%  '

$ws.Range("E46").Value = '  -4.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.26'
$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -2.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("D50").Value = '2.199.66'
$ws.Range("E50").Value = '  -1.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.32'
$ws.Range("E51").Value = '  -1.82%  '
